$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Each of the 15 cells in the 5x3 lattice-multiplication table gets its
# entire contents (problem line, factor line, dashed rule, two partial
# lines) replaced with a new exercise. Table shape (5 rows x 3 cols) is
# unchanged -- only the text inside each cell changes.
$cellData = @(
    @{ Row=1; Col=1; Texts=@("41 x 34", "  3    4", "  ----", "4|    |", "1|    |") },
    @{ Row=1; Col=2; Texts=@("94 x 47", "  4    7", "  ----", "9|    |", "4|    |") },
    @{ Row=1; Col=3; Texts=@("84 x 93", "  9    3", "  ----", "8|    |", "4|    |") },
    @{ Row=2; Col=1; Texts=@("21 x 69", "  6    9", "  ----", "2|    |", "1|    |") },
    @{ Row=2; Col=2; Texts=@("17 x 59", "  5    9", "  ----", "1|    |", "7|    |") },
    @{ Row=2; Col=3; Texts=@("56 x 46", "  4    6", "  ----", "5|    |", "6|    |") },
    @{ Row=3; Col=1; Texts=@("51 x 31", "  3    1", "  ----", "5|    |", "1|    |") },
    @{ Row=3; Col=2; Texts=@("28 x 74", "  7    4", "  ----", "2|    |", "8|    |") },
    @{ Row=3; Col=3; Texts=@("17 x 63", "  6    3", "  ----", "1|    |", "7|    |") },
    @{ Row=4; Col=1; Texts=@("88 x 85", "  8    5", "  ----", "8|    |", "8|    |") },
    @{ Row=4; Col=2; Texts=@("43 x 14", "  1    4", "  ----", "4|    |", "3|    |") },
    @{ Row=4; Col=3; Texts=@("83 x 72", "  7    2", "  ----", "8|    |", "3|    |") },
    @{ Row=5; Col=1; Texts=@("21 x 98", "  9    8", "  ----", "2|    |", "1|    |") },
    @{ Row=5; Col=2; Texts=@("69 x 26", "  2    6", "  ----", "6|    |", "9|    |") },
    @{ Row=5; Col=3; Texts=@("87 x 70", "  7    0", "  ----", "8|    |", "7|    |") }
)

foreach ($item in $cellData) {
    $cell = $t.Cell($item.Row, $item.Col)
    $newText = [string]::Join([char]11, $item.Texts)
    $cell.Range.Text = $newText
}
